$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '26.318.92'
Set-TextValue $ws 'E2' '  -1.11%  '
Set-TextValue $ws 'D3' '1.590.20'
Set-TextValue $ws 'E3' '  -0.33%  '
Set-TextValue $ws 'E4' '  -0.88%  '
Set-TextValue $ws 'D5' '209.88'
Set-TextValue $ws 'E5' '  -0.72%  '
Set-TextValue $ws 'D6' '0.507'
Set-TextValue $ws 'E6' '  -1.03%  '
Set-TextValue $ws 'E7' '  -0.72%  '
Set-TextValue $ws 'D8' '0.0610'
Set-TextValue $ws 'E8' '  -1.13%  '
Set-TextValue $ws 'E9' '  -0.67%  '
Set-TextValue $ws 'E10' '  -0.59%  '
Set-TextValue $ws 'D11' '0.0843'
Set-TextValue $ws 'E11' '  -0.45%  '
Set-TextValue $ws 'D13' '1.588.40'
Set-TextValue $ws 'E13' '  -0.70%  '
Set-TextValue $ws 'D14' '4.07'
Set-TextValue $ws 'E14' '  +0.44%  '
Set-TextValue $ws 'D16' '64.47'
Set-TextValue $ws 'E16' '  -0.42%  '
Set-TextValue $ws 'D17' '26.335.71'
Set-TextValue $ws 'E17' '  -1.16%  '
Set-TextValue $ws 'E18' '  -1.63%  '
Set-TextValue $ws 'D20' '210.86'
Set-TextValue $ws 'E20' '  +1.19%  '
Set-TextValue $ws 'E21' '  -0.62%  '
Set-TextValue $ws 'E22' '  -0.50%  '
Set-TextValue $ws 'E23' '  -4.57%  '
Set-TextValue $ws 'D24' '8.91'
Set-TextValue $ws 'E24' '  -0.78%  '
Set-TextValue $ws 'D25' '144.92'
Set-TextValue $ws 'E25' '  +0.37%  '
Set-TextValue $ws 'E26' '  -0.61%  '
Set-TextValue $ws 'E27' '  -1.23%  '
Set-TextValue $ws 'E28' '  -0.68%  '
Set-TextValue $ws 'D29' '15.27'
Set-TextValue $ws 'E29' '  +0.10%  '
Set-TextValue $ws 'E30' '  -0.60%  '
Set-TextValue $ws 'E31' '  -0.46%  '
Set-TextValue $ws 'D32' '3.20'
Set-TextValue $ws 'E32' '  -1.12%  '
Set-TextValue $ws 'E33' '  +0.41%  '
Set-TextValue $ws 'D34' '1.305.25'
Set-TextValue $ws 'E34' '  +2.31%  '
Set-TextValue $ws 'D35' '0.616'
Set-TextValue $ws 'E35' '  +2.84%  '
Set-TextValue $ws 'E36' '  -2.00%  '
Set-TextValue $ws 'E37' '  -0.52%  '
Set-TextValue $ws 'E38' '  -0.45%  '
Set-TextValue $ws 'E39' '  -13.36%  '
Set-TextValue $ws 'E40' '  -1.42%  '
Set-TextValue $ws 'E41' '  -0.55%  '
Set-TextValue $ws 'E42' '  +3.36%  '
Set-TextValue $ws 'D43' '62.62'
Set-TextValue $ws 'E43' '  +0.17%  '
Set-TextValue $ws 'E44' '  -1.44%  '
Set-TextValue $ws 'E45' '  -1.54%  '
Set-TextValue $ws 'D46' '1.726.14'
Set-TextValue $ws 'E46' '  -0.37%  '
Set-TextValue $ws 'D47' '87.94'
Set-TextValue $ws 'E47' '  -2.56%  '
Set-TextValue $ws 'B48' 'RenderToken'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D48' '1.50'
Set-TextValue $ws 'E48' '  -4.25%  '
Set-TextValue $ws 'B49' 'Algorand'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D49' '0.0983'
Set-TextValue $ws 'E49' '  -4.56%  '
Set-TextValue $ws 'B50' 'Cronos'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D50' '0.0504'
Set-TextValue $ws 'E50' '  -1.57%  '
Set-TextValue $ws 'B51' 'USDD'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue $ws 'D51' '1.00'
Set-TextValue $ws 'E51' '  -0.60%  '
